$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: cell reference, new value. The source values are all plain text
# (inline strings in the sheet XML) even when they look like numbers, e.g. a
# price of "1.784.11" or "1.003". Assigning such a string straight to .Value
# makes Excel auto-detect it as a number when it parses cleanly (e.g. "1.003"),
# which would corrupt values such as "1.003" -> 1.003 (losing the trailing zero)
# or silently drop the thousands-style dot grouping. For those cells we force
# the cell to Text format first, assign the literal string, then restore the
# default "Normal" style so no stray formatting is left behind.
$updates = @(
    ,@("D2", "27.143.42")
    ,@("E2", "  -1.11%  ")
    ,@("D3", "1.784.11")
    ,@("E3", "  -1.72%  ")
    ,@("D4", "1.003")
    ,@("E4", "  +0.13%  ")
    ,@("D5", "336.63")
    ,@("E5", "  -1.95%  ")
    ,@("D6", "1.001")
    ,@("E6", "  +0.14%  ")
    ,@("D7", "0.3831")
    ,@("E7", "  +0.48%  ")
    ,@("D8", "0.3426")
    ,@("E8", "  -2.06%  ")
    ,@("D9", "47.97")
    ,@("D10", "1.191")
    ,@("E10", "  -3.52%  ")
    ,@("D11", "0.07470")
    ,@("E11", "  -3.47%  ")
    ,@("D12", "1.002")
    ,@("E12", "  +0.08%  ")
    ,@("D13", "21.73")
    ,@("E13", "  -1.89%  ")
    ,@("D14", "6.442")
    ,@("E14", "  -2.56%  ")
    ,@("D15", "1.780.03")
    ,@("E15", "  -1.86%  ")
    ,@("D16", "7.116")
    ,@("E16", "  -1.66%  ")
    ,@("D17", "0.00001097")
    ,@("E17", "  -2.14%  ")
    ,@("D18", "0.06650")
    ,@("E18", "  -0.96%  ")
    ,@("D19", "83.45")
    ,@("E19", "  -3.41%  ")
    ,@("E20", "  +0.04%  ")
    ,@("D21", "17.50")
    ,@("E21", "  -0.64%  ")
    ,@("D22", "6.523")
    ,@("E22", "  -0.88%  ")
    ,@("D23", "27.135.97")
    ,@("E23", "  -1.14%  ")
    ,@("D24", "12.29")
    ,@("E24", "  -7.15%  ")
    ,@("E25", "  -3.78%  ")
    ,@("D26", "2.509")
    ,@("E26", "  -6.11%  ")
    ,@("D27", "21.18")
    ,@("E27", "  -3.94%  ")
    ,@("D28", "1.446")
    ,@("E28", "  -1.68%  ")
    ,@("D29", "155.47")
    ,@("E29", "  +0.98%  ")
    ,@("D30", "1.984.44")
    ,@("E30", "  -1.57%  ")
    ,@("D31", "134.31")
    ,@("E31", "  -1.23%  ")
    ,@("D32", "3.978")
    ,@("E32", "  -1.54%  ")
    ,@("D33", "6.035")
    ,@("E33", "  -4.84%  ")
    ,@("D34", "0.08678")
    ,@("E34", "  -1.23%  ")
    ,@("D35", "13.08")
    ,@("E35", "  -6.28%  ")
    ,@("D36", "1.622")
    ,@("E36", "  -4.12%  ")
    ,@("D37", "5.402")
    ,@("E37", "  -4.12%  ")
    ,@("D38", "0.6827")
    ,@("E38", "  -2.13%  ")
    ,@("D39", "0.06336")
    ,@("E39", "  -2.41%  ")
    ,@("D40", "0.02337")
    ,@("E40", "  -2.91%  ")
    ,@("D41", "0.2185")
    ,@("E41", "  -3.96%  ")
    ,@("B42", "TrustWalletToken")
    ,@("C42", "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt")
    ,@("D42", "1.239")
    ,@("E42", "  -4.43%  ")
    ,@("B43", "FraxShare")
    ,@("C43", "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs")
    ,@("D43", "8.435")
    ,@("E43", "  -5.75%  ")
    ,@("D44", "14.26")
    ,@("E44", "  -3.07%  ")
    ,@("D45", "1.000")
    ,@("E45", "  -0.07%  ")
    ,@("D46", "0.6426")
    ,@("E46", "  -1.61%  ")
    ,@("D47", "3.857")
    ,@("E47", "  -3.99%  ")
    ,@("D48", "2.166")
    ,@("E48", "  -0.60%  ")
    ,@("D49", "131.33")
    ,@("E49", "  -1.41%  ")
    ,@("D50", "0.07107")
    ,@("E50", "  -2.98%  ")
    ,@("D51", "78.76")
    ,@("E51", "  -2.35%  ")
)

foreach ($u in $updates) {
    $ref = $u[0]
    $val = $u[1]
    $cell = $ws.Range($ref)
    if ($val -match '^[0-9]+(\.[0-9]+)?$') {
        $cell.NumberFormat = "@"
        $cell.Value = $val
        $cell.Style = "Normal"
    } else {
        $cell.Value = $val
    }
}
